$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 167
$ws.Range("I33").Value = 167
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 167
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = 62
$ws.Range("N33").ClearContents()
$ws.Range("H121").Value = 781.16
$ws.Range("J121").Value = 779.5217
$ws.Range("L121").Value = 2338.5651
$ws.Range("N121").Value = -5832.5651
$ws.Range("H125").Value = 3921.8462
$ws.Range("I125").Value = 1310
$ws.Range("J125").Value = 5554.25
$ws.Range("K125").Value = 11790
$ws.Range("L125").Value = 49988.25
$ws.Range("M125").Value = -9330
$ws.Range("N125").Value = -54908.25
$ws.Range("H132").Value = 2170.3901
$ws.Range("I132").Value = 1765.3438
$ws.Range("J132").Value = 3610.5557
$ws.Range("K132").Value = 5296.0314
$ws.Range("L132").Value = 10831.6671
$ws.Range("M132").Value = -2766.0314
$ws.Range("N132").Value = -15891.6671
$ws.Range("H137").Value = 11112999
$ws.Range("I137").Value = 1670.32
$ws.Range("J137").Value = 25002160
$ws.Range("K137").Value = 5010.96
$ws.Range("L137").Value = 75006480
$ws.Range("M137").Value = -2460.96
$ws.Range("N137").Value = -75011580
$ws.Range("H138").Value = 3318.1516
$ws.Range("I138").Value = 1389.069
$ws.Range("J138").Value = 4117.343
$ws.Range("K138").Value = 4167.207
$ws.Range("L138").Value = 12352.029
$ws.Range("M138").Value = 972.7929999999997
$ws.Range("N138").Value = -22632.029
$ws.Range("H141").Value = 1120.3793
$ws.Range("I141").Value = 1148.44
$ws.Range("J141").Value = 945
$ws.Range("K141").Value = 3445.32
$ws.Range("L141").Value = 2835
$ws.Range("M141").Value = 1734.68
$ws.Range("N141").Value = -13195

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19684.918
$ws.Range("I32").Value = 18172.688
$ws.Range("J32").Value = 32622.889
$ws.Range("K32").Value = 18172.688
$ws.Range("L32").Value = 32622.889
$ws.Range("M32").Value = -17885.688
$ws.Range("N32").Value = -33196.889
$ws.Range("H61").Value = 1704.4667
$ws.Range("I61").Value = 1496
$ws.Range("K61").Value = 1496
$ws.Range("M61").Value = -1284
$ws.Range("H122").Value = 1078.8334
$ws.Range("I122").Value = 1078.8334
$ws.Range("K122").Value = 3236.5002
$ws.Range("M122").Value = -786.5001999999999
$ws.Range("H132").Value = 2228.8
$ws.Range("I132").Value = 1679.6522
$ws.Range("J132").Value = 3281.3333
$ws.Range("K132").Value = 5038.9566
$ws.Range("L132").Value = 9843.999899999999
$ws.Range("M132").Value = -2508.9566
$ws.Range("N132").Value = -14903.9999
$ws.Range("H136").Value = 1704.4667
$ws.Range("I136").Value = 1496
$ws.Range("K136").Value = 4488
$ws.Range("M136").Value = -1938

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 45489520
$ws.Range("I20").Value = 92993.57000000001
$ws.Range("K20").Value = 92993.57000000001
$ws.Range("M20").Value = -92746.57000000001
$ws.Range("H134").Value = 47527
$ws.Range("I134").Value = 2215.1516
$ws.Range("J134").Value = 183462.55
$ws.Range("K134").Value = 6645.4548
$ws.Range("L134").Value = 550387.6499999999
$ws.Range("M134").Value = -4110.4548
$ws.Range("N134").Value = -555457.6499999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1784.1305
$ws.Range("I31").Value = 777.8461
$ws.Range("J31").Value = 2180.5454
$ws.Range("K31").Value = 777.8461
$ws.Range("L31").Value = 2180.5454
$ws.Range("M31").Value = -482.8461
$ws.Range("N31").Value = -2770.5454
$ws.Range("H34").Value = 1784.1305
$ws.Range("I34").Value = 777.8461
$ws.Range("J34").Value = 2180.5454
$ws.Range("K34").Value = 777.8461
$ws.Range("L34").Value = 2180.5454
$ws.Range("M34").Value = -575.8461
$ws.Range("N34").Value = -2584.5454
$ws.Range("H105").Value = 2058.1853
$ws.Range("I105").Value = 1738.55
$ws.Range("J105").Value = 2971.4285
$ws.Range("K105").Value = 1738.55
$ws.Range("L105").Value = 2971.4285
$ws.Range("M105").Value = 8.450000000000045
$ws.Range("N105").Value = -6465.4285

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4509.0454
$ws.Range("I70").Value = 4094.1177
$ws.Range("K70").Value = 4094.1177
$ws.Range("M70").Value = -3824.1177
$ws.Range("H73").Value = 4509.0454
$ws.Range("I73").Value = 4094.1177
$ws.Range("K73").Value = 4094.1177
$ws.Range("M73").Value = -3158.1177
$ws.Range("H126").Value = 3285.5715
$ws.Range("I126").Value = 2999.8
$ws.Range("J126").Value = 4000
$ws.Range("K126").Value = 8999.400000000001
$ws.Range("L126").Value = 12000
$ws.Range("M126").Value = -6529.400000000001
$ws.Range("N126").Value = -16940
$ws.Range("H132").Value = 2398.718
$ws.Range("I132").Value = 1298.174
$ws.Range("K132").Value = 3894.522
$ws.Range("M132").Value = -1364.522

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1224.625
$ws.Range("I7").Value = 1179.4
$ws.Range("J7").Value = 1300
$ws.Range("K7").Value = 1179.4
$ws.Range("L7").Value = 1300
$ws.Range("M7").Value = -1067.4
$ws.Range("N7").Value = -1524
$ws.Range("H122").Value = 5779.696
$ws.Range("I122").Value = 6249.4
$ws.Range("K122").Value = 18748.2
$ws.Range("M122").Value = -16298.2
$ws.Range("H126").Value = 1224.625
$ws.Range("I126").Value = 1179.4
$ws.Range("J126").Value = 1300
$ws.Range("K126").Value = 3538.2
$ws.Range("L126").Value = 3900
$ws.Range("M126").Value = -1068.2
$ws.Range("N126").Value = -8840
$ws.Range("H132").Value = 2022770
$ws.Range("I132").Value = 2599761.5
$ws.Range("J132").Value = 3299.8
$ws.Range("K132").Value = 7799284.5
$ws.Range("L132").Value = 9899.400000000001
$ws.Range("M132").Value = -7796754.5
$ws.Range("N132").Value = -14959.4
$ws.Range("H136").Value = 1918.9678
$ws.Range("I136").Value = 1159.921
$ws.Range("J136").Value = 3120.7917
$ws.Range("K136").Value = 3479.763
$ws.Range("L136").Value = 9362.375100000001
$ws.Range("M136").Value = -929.7629999999999
$ws.Range("N136").Value = -14462.3751

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1674
$ws.Range("I122").Value = 1501.5
$ws.Range("J122").Value = 1881
$ws.Range("K122").Value = 4504.5
$ws.Range("L122").Value = 5643
$ws.Range("M122").Value = -2054.5
$ws.Range("N122").Value = -10543
$ws.Range("H126").Value = 974
$ws.Range("I126").Value = 385.6
$ws.Range("J126").Value = 2445
$ws.Range("K126").Value = 1156.8
$ws.Range("L126").Value = 7335
$ws.Range("M126").Value = 1313.2
$ws.Range("N126").Value = -12275
